$wb = $excel.ActiveWorkbook

# This script applies updated FFXIV "Ixion" market-board price/profit
# figures (columns H-N) produced by the scheduled market-data runner.
# Values were pulled fresh from the market API, so this is a plain data
# refresh -- no formulas, formatting, or structure changes.

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1108.9656
$ws.Range("I112").Value = 1200
$ws.Range("J112").Value = 1102.2222
$ws.Range("K112").Value = 3600
$ws.Range("L112").Value = 3306.6666
$ws.Range("M112").Value = -2492
$ws.Range("N112").Value = -5522.6666
$ws.Range("H113").Value = 25002460
$ws.Range("I113").Value = 3075
$ws.Range("J113").Value = 125000000
$ws.Range("K113").Value = 3075
$ws.Range("L113").Value = 125000000
$ws.Range("M113").Value = 179
$ws.Range("N113").Value = -125006508
$ws.Range("H129").Value = 953.7538500000001
$ws.Range("J129").Value = 979.8246
$ws.Range("L129").Value = 2939.4738
$ws.Range("N129").Value = -12939.4738
$ws.Range("H132").Value = 764.7820400000001
$ws.Range("I132").Value = 563.9559
$ws.Range("J132").Value = 2130.4
$ws.Range("K132").Value = 1691.8677
$ws.Range("L132").Value = 6391.200000000001
$ws.Range("M132").Value = 838.1322999999998
$ws.Range("N132").Value = -11451.2
$ws.Range("H141").Value = 1354.3529
$ws.Range("I141").Value = 926.1707
$ws.Range("J141").Value = 3109.9
$ws.Range("K141").Value = 2778.5121
$ws.Range("L141").Value = 9329.700000000001
$ws.Range("M141").Value = 2401.4879
$ws.Range("N141").Value = -19689.7

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2482.209
$ws.Range("I61").Value = 2629.7144
$ws.Range("J61").Value = 1731.2727
$ws.Range("K61").Value = 2629.7144
$ws.Range("L61").Value = 1731.2727
$ws.Range("M61").Value = -2417.7144
$ws.Range("N61").Value = -2155.2727
$ws.Range("H74").Value = 1308.1875
$ws.Range("I74").Value = 1172.2084
$ws.Range("J74").Value = 1716.125
$ws.Range("K74").Value = 1172.2084
$ws.Range("L74").Value = 1716.125
$ws.Range("M74").Value = -298.2084
$ws.Range("N74").Value = -3464.125
$ws.Range("H77").Value = 1308.1875
$ws.Range("I77").Value = 1172.2084
$ws.Range("J77").Value = 1716.125
$ws.Range("K77").Value = 5861.041999999999
$ws.Range("L77").Value = 8580.625
$ws.Range("M77").Value = -1493.041999999999
$ws.Range("N77").Value = -17316.625
$ws.Range("H97").Value = 1084.7894
$ws.Range("I97").Value = 899.9231
$ws.Range("J97").Value = 1485.3334
$ws.Range("K97").Value = 899.9231
$ws.Range("L97").Value = 1485.3334
$ws.Range("M97").Value = -403.9231
$ws.Range("N97").Value = -2477.3334
$ws.Range("H122").Value = 1510905.8
$ws.Range("I122").Value = 2139949.8
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 6419849.399999999
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -6417399.399999999
$ws.Range("N122").Value = -8500
$ws.Range("H132").Value = 1756772.4
$ws.Range("I132").Value = 1756.3414
$ws.Range("J132").Value = 6254001
$ws.Range("K132").Value = 5269.0242
$ws.Range("L132").Value = 18762003
$ws.Range("M132").Value = -2739.0242
$ws.Range("N132").Value = -18767063
$ws.Range("H136").Value = 2482.209
$ws.Range("I136").Value = 2629.7144
$ws.Range("J136").Value = 1731.2727
$ws.Range("K136").Value = 7889.1432
$ws.Range("L136").Value = 5193.8181
$ws.Range("M136").Value = -5339.1432
$ws.Range("N136").Value = -10293.8181

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4711.8237
$ws.Range("I134").Value = 5080.4287
$ws.Range("J134").Value = 2991.6667
$ws.Range("K134").Value = 15241.2861
$ws.Range("L134").Value = 8975.000100000001
$ws.Range("M134").Value = -12706.2861
$ws.Range("N134").Value = -14045.0001

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 100000
$ws.Range("I2").Value = 100000
$ws.Range("K2").Value = 100000
$ws.Range("M2").Value = -99887
$ws.Range("H31").Value = 5628.3716
$ws.Range("I31").Value = 1543.6511
$ws.Range("J31").Value = 10646.743
$ws.Range("K31").Value = 1543.6511
$ws.Range("L31").Value = 10646.743
$ws.Range("M31").Value = -1248.6511
$ws.Range("N31").Value = -11236.743
$ws.Range("H34").Value = 5628.3716
$ws.Range("I34").Value = 1543.6511
$ws.Range("J34").Value = 10646.743
$ws.Range("K34").Value = 1543.6511
$ws.Range("L34").Value = 10646.743
$ws.Range("M34").Value = -1341.6511
$ws.Range("N34").Value = -11050.743
$ws.Range("H58").Value = 1210.4
$ws.Range("I58").Value = 775.7368
$ws.Range("J58").Value = 2586.8333
$ws.Range("K58").Value = 775.7368
$ws.Range("L58").Value = 2586.8333
$ws.Range("M58").Value = -572.7368
$ws.Range("N58").Value = -2992.8333
$ws.Range("H132").Value = 1744.6792
$ws.Range("I132").Value = 1378.5278
$ws.Range("K132").Value = 4135.5834
$ws.Range("M132").Value = -1605.5834
$ws.Range("H136").Value = 1210.4
$ws.Range("I136").Value = 775.7368
$ws.Range("J136").Value = 2586.8333
$ws.Range("K136").Value = 2327.2104
$ws.Range("L136").Value = 7760.499899999999
$ws.Range("M136").Value = 222.7896000000001
$ws.Range("N136").Value = -12860.4999

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2047.7778
$ws.Range("J39").Value = 2047.7778
$ws.Range("L39").Value = 6143.3334
$ws.Range("N39").Value = -6731.3334
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1766.0577
$ws.Range("I132").Value = 1466.8462
$ws.Range("J132").Value = 2663.6924
$ws.Range("K132").Value = 4400.5386
$ws.Range("L132").Value = 7991.0772
$ws.Range("M132").Value = -1870.5386
$ws.Range("N132").Value = -13051.0772

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5666.841
$ws.Range("I136").Value = 4855.1714
$ws.Range("J136").Value = 8823.333000000001
$ws.Range("K136").Value = 14565.5142
$ws.Range("L136").Value = 26469.999
$ws.Range("M136").Value = -12015.5142
$ws.Range("N136").Value = -31569.999

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 6035
$ws.Range("I51").Value = 6035
$ws.Range("K51").Value = 6035
$ws.Range("M51").Value = -5525
$ws.Range("H132").Value = 17457.383
$ws.Range("I132").Value = 23339.455
$ws.Range("J132").Value = 1281.6875
$ws.Range("K132").Value = 70018.36500000001
$ws.Range("L132").Value = 3845.0625
$ws.Range("M132").Value = -67488.36500000001
$ws.Range("N132").Value = -8905.0625
$ws.Range("H136").Value = 5496207
$ws.Range("I136").Value = 1694.4849
$ws.Range("K136").Value = 5083.4547
$ws.Range("M136").Value = -2533.4547
